$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/19/2024  Through  2/25/2024"

# --- Row 14 ---
# M14 changes from text placeholder "***.*" to a real number; pull the
# numeric style (s=16) from L14 (same row, untouched) before writing the value.
$ws.Range("L14").Copy($ws.Range("M14"))
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -88.888888888888

# --- Row 15 ---
# C15 and D15 flip from numbers to the text placeholder "0" (style s=14);
# E15 flips from a number to the text placeholder "***.*" (style s=14).
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("H14").Copy($ws.Range("E15"))
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("L15").Value = 25
$ws.Range("N15").Value = -44.444444444444

# --- Row 16 ---
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 175
$ws.Range("F16").Value = 43
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = 26.470588235294
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 66
$ws.Range("K16").Value = 34.848484848484
$ws.Range("L16").Value = 34.848484848484
$ws.Range("M16").Value = 34.848484848484
$ws.Range("N16").Value = -67.037037037037

# --- Row 17 ---
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -35.294117647058
$ws.Range("F17").Value = 70
$ws.Range("G17").Value = 70
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = 8.130081300813
$ws.Range("L17").Value = 34.343434343434
$ws.Range("M17").Value = 125.423728813559
$ws.Range("N17").Value = 24.299065420560

# --- Row 18 ---
$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 180
$ws.Range("F18").Value = 31
$ws.Range("H18").Value = 72.222222222222
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 35.897435897435
$ws.Range("L18").Value = -14.516129032258
$ws.Range("M18").Value = 65.625
$ws.Range("N18").Value = -74.641148325358

# --- Row 19 ---
$ws.Range("C19").Value = 20
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 65
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 47.727272727272
$ws.Range("I19").Value = 136
$ws.Range("J19").Value = 103
$ws.Range("K19").Value = 32.038834951456
$ws.Range("L19").Value = 43.157894736842
$ws.Range("M19").Value = 138.59649122807
$ws.Range("N19").Value = 51.111111111111

# --- Row 20 ---
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = -76.470588235294
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -37.931034482758
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -42.105263157894
$ws.Range("L20").Value = -41.071428571428
$ws.Range("M20").Value = 83.333333333333
$ws.Range("N20").Value = -65.625

# --- Row 21 ---
$ws.Range("C21").Value = 60
$ws.Range("D21").Value = 53
$ws.Range("E21").Value = 13.207547169811
$ws.Range("F21").Value = 230
$ws.Range("G21").Value = 198
$ws.Range("H21").Value = 16.161616161616
$ws.Range("I21").Value = 451
$ws.Range("J21").Value = 398
$ws.Range("K21").Value = 13.316582914572
$ws.Range("L21").Value = 17.447916666666
$ws.Range("M21").Value = 91.101694915254
$ws.Range("N21").Value = -43.554443053817

# --- Row 22 ---
$ws.Range("C22").Value = 3
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("H14").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 17
$ws.Range("K22").Value = 142.857142857143
$ws.Range("L22").Value = -22.727272727272
$ws.Range("M22").Value = 112.5

# --- Row 23 ---
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 33
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = 3.125
$ws.Range("I23").Value = 69
$ws.Range("J23").Value = 63
$ws.Range("K23").Value = 9.523809523809
$ws.Range("L23").Value = 43.75
$ws.Range("M23").Value = 60.465116279069

# --- Row 24 ---
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 65
$ws.Range("F24").Value = 130
$ws.Range("G24").Value = 103
$ws.Range("H24").Value = 26.213592233009
$ws.Range("I24").Value = 274
$ws.Range("J24").Value = 238
$ws.Range("K24").Value = 15.126050420168
$ws.Range("L24").Value = 28.638497652582
$ws.Range("M24").Value = 50.549450549450

# --- Row 25 ---
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = -14.814814814814
$ws.Range("F25").Value = 103
$ws.Range("G25").Value = 77
$ws.Range("H25").Value = 33.766233766233
$ws.Range("I25").Value = 196
$ws.Range("J25").Value = 147
$ws.Range("K25").Value = 33.333333333333
$ws.Range("L25").Value = 53.125
$ws.Range("M25").Value = 36.111111111111

# --- Row 26 ---
$ws.Range("C14").Copy($ws.Range("C26"))
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -57.142857142857
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = -50
$ws.Range("L26").Value = 16.666666666666

# --- Row 27 ---
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = -36.363636363636
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -17.647058823529
$ws.Range("L27").Value = 27.272727272727

# --- Row 28 ---
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("F28").Value = 3
$ws.Range("M28").Value = -73.333333333333
$ws.Range("N28").Value = -87.878787878787

# --- Row 29 ---
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("F29").Value = 3
$ws.Range("M29").Value = -55.555555555555
$ws.Range("N29").Value = -84.615384615384
